{"js": "const pairs = [\n  [\"797\u00d77=5579\", \"403\u00d74=1612\"],\n  [\"822\u00d76=4932\", \"812\u00d78=6496\"],\n  [\"654\u00d75=3270\", \"146\u00d79=1314\"],\n  [\"285\u00d79=2565\", \"620\u00d75=3100\"],\n  [\"544\u00d72=1088\", \"403\u00d77=2821\"],\n  [\"312\u00d72=624\", \"623\u00d74=2492\"],\n  [\"349\u00d78=2792\", \"217\u00d76=1302\"],\n  [\"710\u00d72=1420\", \"303\u00d78=2424\"],\n  [\"379\u00d73=1137\", \"912\u00d78=7296\"],\n  [\"683\u00d77=4781\", \"836\u00d73=2508\"],\n  [\"148\u00d77=1036\", \"105\u00d74=420\"],\n  [\"646\u00d76=3876\", \"767\u00d72=1534\"],\n  [\"976\u00d79=8784\", \"440\u00d72=880\"],\n  [\"787\u00d74=3148\", \"754\u00d79=6786\"],\n  [\"708\u00d74=2832\", \"185\u00d72=370\"],\n  [\"142\u00d79=1278\", \"841\u00d75=4205\"],\n  [\"246\u00d76=1476\", \"206\u00d79=1854\"],\n  [\"247\u00d77=1729\", \"574\u00d75=2870\"],\n  [\"283\u00d73=849\", \"780\u00d74=3120\"],\n  [\"494\u00d79=4446\", \"580\u00d73=1740\"],\n  [\"678\u00d74=2712\", \"306\u00d78=2448\"],\n  [\"241\u00d73=723\", \"965\u00d73=2895\"],\n  [\"591\u00d74=2364\", \"976\u00d73=2928\"],\n  [\"681\u00d78=5448\", \"854\u00d73=2562\"],\n  [\"696\u00d76=4176\", \"992\u00d72=1984\"],\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\nfor (const [findText, replaceText] of pairs) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n    totalReplaced++;\n  }\n  await context.sync();\n}\nreturn totalReplaced;\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"797\u00d77=5579\", \"403\u00d74=1612\")\n  ,@(\"822\u00d76=4932\", \"812\u00d78=6496\")\n  ,@(\"654\u00d75=3270\", \"146\u00d79=1314\")\n  ,@(\"285\u00d79=2565\", \"620\u00d75=3100\")\n  ,@(\"544\u00d72=1088\", \"403\u00d77=2821\")\n  ,@(\"312\u00d72=624\", \"623\u00d74=2492\")\n  ,@(\"349\u00d78=2792\", \"217\u00d76=1302\")\n  ,@(\"710\u00d72=1420\", \"303\u00d78=2424\")\n  ,@(\"379\u00d73=1137\", \"912\u00d78=7296\")\n  ,@(\"683\u00d77=4781\", \"836\u00d73=2508\")\n  ,@(\"148\u00d77=1036\", \"105\u00d74=420\")\n  ,@(\"646\u00d76=3876\", \"767\u00d72=1534\")\n  ,@(\"976\u00d79=8784\", \"440\u00d72=880\")\n  ,@(\"787\u00d74=3148\", \"754\u00d79=6786\")\n  ,@(\"708\u00d74=2832\", \"185\u00d72=370\")\n  ,@(\"142\u00d79=1278\", \"841\u00d75=4205\")\n  ,@(\"246\u00d76=1476\", \"206\u00d79=1854\")\n  ,@(\"247\u00d77=1729\", \"574\u00d75=2870\")\n  ,@(\"283\u00d73=849\", \"780\u00d74=3120\")\n  ,@(\"494\u00d79=4446\", \"580\u00d73=1740\")\n  ,@(\"678\u00d74=2712\", \"306\u00d78=2448\")\n  ,@(\"241\u00d73=723\", \"965\u00d73=2895\")\n  ,@(\"591\u00d74=2364\", \"976\u00d73=2928\")\n  ,@(\"681\u00d78=5448\", \"854\u00d73=2562\")\n  ,@(\"696\u00d76=4176\", \"992\u00d72=1984\")\n)\n\n$replacedCount = 0\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  $ok = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n  if ($ok) { $replacedCount = $replacedCount + 1 }\n}\n\n\"Replaced: $replacedCount\"\n"}
